# Update countries & provincias Spain
# Applies the 16-Aug-2020 11:01 data refresh: updated case counters for a
# number of countries/regions, the refreshed timestamp footer, and the
# consequent re-ranking (by total cases) that swaps a few adjacent rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 16 de Agosto de 2020 a las 11:01"

# --- Re-ranked adjacent country pairs (names swap places) -------------
# Rusia keeps its rank (row 6) but gets refreshed counters.
$ws.Range("B6").Value = 2594112
$ws.Range("C6").Value = 4904
$ws.Range("D6").Value = 1862937
$ws.Range("E6").Value = 681053
$ws.Range("G6").Value = 38
$ws.Range("H6").Value = 50122

$ws.Range("B26").Value = 139549
$ws.Range("C26").Value = 2081
$ws.Range("D26").Value = 93103
$ws.Range("E26").Value = 40296
$ws.Range("G26").Value = 79
$ws.Range("H26").Value = 6150

$ws.Range("B33").Value = 92404
$ws.Range("C33").Value = 171
$ws.Range("D33").Value = 68400
$ws.Range("E33").Value = 23325
$ws.Range("G33").Value = 5
$ws.Range("H33").Value = 679

$ws.Range("B34").Value = 91356
$ws.Range("C34").Value = 1637
$ws.Range("D34").Value = 47822
$ws.Range("E34").Value = 41464
$ws.Range("G34").Value = 26
$ws.Range("H34").Value = 2070

$ws.Range("B47").Value = 56684
$ws.Range("C47").Value = 594
$ws.Range("D47").Value = 39130
$ws.Range("E47").Value = 15677
$ws.Range("G47").Value = 8
$ws.Range("H47").Value = 1877

$ws.Range("B48").Value = 55747
$ws.Range("C48").Value = 86
$ws.Range("E48").Value = 4199

# Australia (rank 75) <-> Austria (rank 76) swap places.
$ws.Range("A71").Value = "Austria"
$ws.Range("B71").Value = 23370
$ws.Range("C71").Value = 191
$ws.Range("D71").Value = 20681
$ws.Range("E71").Value = 1961
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 728

$ws.Range("A72").Value = "Australia"
$ws.Range("B72").Value = 23287
$ws.Range("C72").Value = 252
$ws.Range("D72").Value = 13634
$ws.Range("E72").Value = 9257
$ws.Range("G72").Value = 17
$ws.Range("H72").Value = 396

# Malasia (rank 93) <-> Zambia (rank 94) swap places.
$ws.Range("A89").Value = "Malasia"
$ws.Range("B89").Value = 9200
$ws.Range("C89").Value = 25
$ws.Range("D89").Value = 8859
$ws.Range("E89").Value = 216
$ws.Range("H89").Value = 125

$ws.Range("A90").Value = "Zambia"
$ws.Range("B90").Value = 9186
$ws.Range("D90").Value = 8065
$ws.Range("E90").Value = 861
$ws.Range("H90").Value = 260

$ws.Range("B111").Value = 4481
$ws.Range("C111").Value = 74
$ws.Range("D111").Value = 3548
$ws.Range("E111").Value = 864
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 69

# Eslovaquia (rank 127) <-> Sri Lanka (rank 128) swap places.
$ws.Range("A123").Value = "Eslovaquia"
$ws.Range("B123").Value = 2902
$ws.Range("C123").Value = 47
$ws.Range("D123").Value = 1969
$ws.Range("E123").Value = 902
$ws.Range("H123").Value = 31

$ws.Range("A124").Value = "Sri Lanka"
$ws.Range("B124").Value = 2890
$ws.Range("D124").Value = 2666
$ws.Range("E124").Value = 213
$ws.Range("H124").Value = 11

$ws.Range("B128").Value = 2416
$ws.Range("C128").Value = 15
$ws.Range("D128").Value = 2051
$ws.Range("E128").Value = 236

$ws.Range("B131").Value = 2190
$ws.Range("C131").Value = 6
$ws.Range("E131").Value = 151

$ws.Range("B168").Value = 484
$ws.Range("C168").Value = 2
$ws.Range("E168").Value = 27

# Islas Malvinas (rank 217) <-> Montserrat (rank 218) swap places.
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
